$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Shift_coded (column B) values for the specified rows
$ws.Range("B4").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("B13").Value = 7
$ws.Range("B16").Value = 3
$ws.Range("B18").Value = 3
$ws.Range("B20").Value = 3
$ws.Range("B21").Value = 6
$ws.Range("B22").Value = 7
$ws.Range("B27").Value = 7
$ws.Range("B30").Value = 4
$ws.Range("B32").Value = 4
$ws.Range("B34").Value = 4
$ws.Range("B36").Value = 6
$ws.Range("B45").Value = 5
$ws.Range("B47").Value = 5
$ws.Range("B48").Value = 5
$ws.Range("B49").Value = 6

# Remove the AutoFilter from the sheet
$ws.AutoFilterMode = $false

# Select cell B1 so it becomes the active cell / selection in the sheet view
$ws.Range("B1").Select()
